$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: numeric-looking price strings must be written with a leading
# apostrophe so Excel stores them as literal text (matching the workbook's
# existing inlineStr/text convention) instead of converting them to numbers.

# --- Simple Price (column D) updates ---
$ws.Range("D2").Value  = "'264.79"
$ws.Range("D3").Value  = "'22.78"
$ws.Range("D4").Value  = "'6.232"
$ws.Range("D5").Value  = "'0.06118"
$ws.Range("D6").Value  = "'3.557"
$ws.Range("D7").Value  = "'6.730"
$ws.Range("D8").Value  = "'1.377"
$ws.Range("D9").Value  = "'0.8135"
$ws.Range("D10").Value = "'0.1593"
$ws.Range("D11").Value = "'0.08218"
$ws.Range("D13").Value = "'0.03178"
$ws.Range("D14").Value = "'0.09251"
$ws.Range("D15").Value = "'3.906"
$ws.Range("D16").Value = "'0.001686"
$ws.Range("D17").Value = "'0.04851"
$ws.Range("D18").Value = "'0.0006277"
$ws.Range("D19").Value = "'0.006237"
$ws.Range("D20").Value = "'0.001105"
$ws.Range("D21").Value = "'0.003201"
$ws.Range("D22").Value = "'0.0001505"
$ws.Range("D23").Value = "'3.692"
$ws.Range("D26").Value = "'0.1273"
$ws.Range("D27").Value = "'0.0002693"
$ws.Range("D40").Value = "'0.04588"

# --- Rows 41-43: coins rotated (KickToken, BKEXToken, CEJI) with new prices ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007285"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1124"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003400"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining Price (column D) updates ---
$ws.Range("D45").Value = "'0.00006177"
$ws.Range("D47").Value = "'0.7532"
$ws.Range("D48").Value = "'0.2383"
$ws.Range("D49").Value = "'0.00002108"
